$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 (right after the header row).
# This shifts all existing data rows down by one, which matches the
# diff exactly: old row N's content becomes row N+1's content, and the
# old last row (25) ends up surviving as the new row 26.
$ws.Rows.Item(2).Insert()

# The Insert() operation copies the formatting of the row above (the
# bold/centered header style) into the freshly inserted row. Reset the
# new row's style back to the plain "no style" look used by all the
# other data rows before we populate it.
$ws.Range("A2:R2").ClearFormats()

# Populate the new row 2 with this week's new record.
$ws.Range("A2").Value = 5
$ws.Range("B2").Value = "Macroferia Regional de Talca"
$ws.Range("C2").Value = "Maule"
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D2").Value = 44630
$ws.Range("E2").Value = 7
$ws.Range("F2").Value = 100112043
$ws.Range("G2").Value = "Pepino dulce"
$ws.Range("H2").Value = "Cultivar IV Región"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 300
$ws.Range("K2").Value = 15000
$ws.Range("L2").Value = 15000
$ws.Range("M2").Value = 15000
$ws.Range("N2").Value = "$/bandeja 18 kilos"
$ws.Range("O2").Value = "Provincia de Limarí"
$ws.Range("P2").Value = 833
$ws.Range("Q2").Value = 18
$ws.Range("R2").Value = "Hortaliza"
